$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 20,10
$data[0,0] = -19.96870897467064
$data[0,1] = 2.628207110207388
$data[0,2] = -19.96870897467064
$data[0,3] = -19.96870897467064
$data[0,4] = -19.96870897467064
$data[0,5] = -19.96870897467064
$data[0,6] = -19.96870897467064
$data[0,7] = -19.96870897467064
$data[0,8] = -19.96870897467064
$data[0,9] = -19.96870897467064
$data[1,0] = -19.96870897467064
$data[1,1] = -19.96870897467064
$data[1,2] = -19.96870897467064
$data[1,3] = -19.96870897467064
$data[1,4] = -19.96870897467064
$data[1,5] = -19.96870897467064
$data[1,6] = -19.96870897467064
$data[1,7] = 2.321166107383742
$data[1,8] = -19.96870897467064
$data[1,9] = -19.96870897467064
$data[2,0] = -19.96870897467064
$data[2,1] = 2.129716731662368
$data[2,2] = 2.895794020233535
$data[2,3] = -19.96870897467064
$data[2,4] = 2.461258809566237
$data[2,5] = -19.96870897467064
$data[2,6] = 1.795758113100816
$data[2,7] = -19.96870897467064
$data[2,8] = 2.355534486713058
$data[2,9] = -19.96870897467064
$data[3,0] = -19.96870897467064
$data[3,1] = 0.4852779813189747
$data[3,2] = -19.96870897467064
$data[3,3] = -19.96870897467064
$data[3,4] = -19.96870897467064
$data[3,5] = 2.108824256136363
$data[3,6] = -19.96870897467064
$data[3,7] = -19.96870897467064
$data[3,8] = -19.96870897467064
$data[3,9] = -19.96870897467064
$data[4,0] = -19.96870897467064
$data[4,1] = -19.96870897467064
$data[4,2] = -19.96870897467064
$data[4,3] = -19.96870897467064
$data[4,4] = -19.96870897467064
$data[4,5] = -19.96870897467064
$data[4,6] = -19.96870897467064
$data[4,7] = -19.96870897467064
$data[4,8] = -19.96870897467064
$data[4,9] = -19.96870897467064
$data[5,0] = 3.003506763826863
$data[5,1] = -19.96870897467064
$data[5,2] = -19.96870897467064
$data[5,3] = -19.96870897467064
$data[5,4] = -19.96870897467064
$data[5,5] = -19.96870897467064
$data[5,6] = -19.96870897467064
$data[5,7] = -19.96870897467064
$data[5,8] = -19.96870897467064
$data[5,9] = -19.96870897467064
$data[6,0] = -19.96870897467064
$data[6,1] = -19.96870897467064
$data[6,2] = -19.96870897467064
$data[6,3] = 4.321926759159549
$data[6,4] = -19.96870897467064
$data[6,5] = -19.96870897467064
$data[6,6] = -19.96870897467064
$data[6,7] = -19.96870897467064
$data[6,8] = -19.96870897467064
$data[6,9] = -19.96870897467064
$data[7,0] = 3.582617801225822
$data[7,1] = -19.96870897467064
$data[7,2] = -19.96870897467064
$data[7,3] = -19.96870897467064
$data[7,4] = -19.96870897467064
$data[7,5] = -19.96870897467064
$data[7,6] = -19.96870897467064
$data[7,7] = -19.96870897467064
$data[7,8] = -19.96870897467064
$data[7,9] = -19.96870897467064
$data[8,0] = -19.96870897467064
$data[8,1] = -19.96870897467064
$data[8,2] = -19.96870897467064
$data[8,3] = -19.96870897467064
$data[8,4] = -19.96870897467064
$data[8,5] = -19.96870897467064
$data[8,6] = -19.96870897467064
$data[8,7] = 1.525559212052394
$data[8,8] = -19.96870897467064
$data[8,9] = 2.210975709162597
$data[9,0] = -19.96870897467064
$data[9,1] = -19.96870897467064
$data[9,2] = -19.96870897467064
$data[9,3] = -19.96870897467064
$data[9,4] = -19.96870897467064
$data[9,5] = 2.581948813595798
$data[9,6] = -19.96870897467064
$data[9,7] = -19.96870897467064
$data[9,8] = -19.96870897467064
$data[9,9] = 1.354176951929801
$data[10,0] = -19.96870897467064
$data[10,1] = -19.96870897467064
$data[10,2] = -19.96870897467064
$data[10,3] = -19.96870897467064
$data[10,4] = -19.96870897467064
$data[10,5] = -19.96870897467064
$data[10,6] = -19.96870897467064
$data[10,7] = -19.96870897467064
$data[10,8] = -19.96870897467064
$data[10,9] = -19.96870897467064
$data[11,0] = -19.96870897467064
$data[11,1] = -19.96870897467064
$data[11,2] = -19.96870897467064
$data[11,3] = -19.96870897467064
$data[11,4] = -19.96870897467064
$data[11,5] = -19.96870897467064
$data[11,6] = -19.96870897467064
$data[11,7] = -19.96870897467064
$data[11,8] = 2.29622636369992
$data[11,9] = 1.610966505275419
$data[12,0] = -19.96870897467064
$data[12,1] = -19.96870897467064
$data[12,2] = 1.609188636117485
$data[12,3] = -19.96870897467064
$data[12,4] = -19.96870897467064
$data[12,5] = -19.96870897467064
$data[12,6] = -19.96870897467064
$data[12,7] = -19.96870897467064
$data[12,8] = -19.96870897467064
$data[12,9] = 2.132047350888165
$data[13,0] = -19.96870897467064
$data[13,1] = -19.96870897467064
$data[13,2] = -0.1159767955767304
$data[13,3] = -19.96870897467064
$data[13,4] = -19.96870897467064
$data[13,5] = -19.96870897467064
$data[13,6] = -19.96870897467064
$data[13,7] = -19.96870897467064
$data[13,8] = -19.96870897467064
$data[13,9] = -19.96870897467064
$data[14,0] = -19.96870897467064
$data[14,1] = -19.96870897467064
$data[14,2] = -19.96870897467064
$data[14,3] = -19.96870897467064
$data[14,4] = -19.96870897467064
$data[14,5] = -19.96870897467064
$data[14,6] = -19.96870897467064
$data[14,7] = -19.96870897467064
$data[14,8] = 2.319741751789172
$data[14,9] = -19.96870897467064
$data[15,0] = -19.96870897467064
$data[15,1] = -0.07056080327256271
$data[15,2] = -0.8629615632208952
$data[15,3] = -19.96870897467064
$data[15,4] = -19.96870897467064
$data[15,5] = -19.96870897467064
$data[15,6] = 0.5415964370459784
$data[15,7] = 0.8404947059925117
$data[15,8] = 1.276883769270588
$data[15,9] = -19.96870897467064
$data[16,0] = -19.96870897467064
$data[16,1] = -19.96870897467064
$data[16,2] = -19.96870897467064
$data[16,3] = -19.96870897467064
$data[16,4] = -19.96870897467064
$data[16,5] = -19.96870897467064
$data[16,6] = 0.4225860175729358
$data[16,7] = 0.9135527008564504
$data[16,8] = 1.353241684164059
$data[16,9] = -19.96870897467064
$data[17,0] = -19.96870897467064
$data[17,1] = -19.96870897467064
$data[17,2] = 1.625019966647469
$data[17,3] = -19.96870897467064
$data[17,4] = -19.96870897467064
$data[17,5] = -19.96870897467064
$data[17,6] = 1.914104578126343
$data[17,7] = 2.12839887631533
$data[17,8] = -19.96870897467064
$data[17,9] = -19.96870897467064
$data[18,0] = -19.96870897467064
$data[18,1] = 1.774873486041458
$data[18,2] = 2.307313453492803
$data[18,3] = -19.96870897467064
$data[18,4] = 3.857285570632366
$data[18,5] = -19.96870897467064
$data[18,6] = 2.243495570876016
$data[18,7] = 2.027501402242016
$data[18,8] = -19.96870897467064
$data[18,9] = 2.426435233774726
$data[19,0] = -19.96870897467064
$data[19,1] = 1.874559228754789
$data[19,2] = -19.96870897467064
$data[19,3] = -19.96870897467064
$data[19,4] = -19.96870897467064
$data[19,5] = 3.277852023688649
$data[19,6] = 2.386222891603272
$data[19,7] = -19.96870897467064
$data[19,8] = -19.96870897467064
$data[19,9] = -19.96870897467064
$ws.Range("B2:K21").Value = $data
